# Connect to the running Excel instance / active workbook.
$excel = New-Object -ComObject Excel.Application
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: the trade's close time moved later in the day (sped-up trading).
$ws.Cells.Item(4, 7).Value = 42606.498645833337

# Three new trade rows (5-7) appended beneath the existing data.
$tradeRows = @(
    @{ A = 9820.09;             B = 9875.39;  C = 316.81; D = 318.58999999999997; E = $true; F = 0.56000000000000005; G = 42606.585578703707; H = $false },
    @{ A = 9752.33;              B = 9820.09;  C = 316.81; D = 319;                 E = $true; F = 0.69;                 G = 42606.586689814816; H = $false },
    @{ A = 9697.7199999999993;   B = 9752.33;  C = 316.81; D = 318.58999999999997; E = $true; F = 0.56000000000000005; G = 42606.587800925925; H = $false }
)

$startRow = 5
$endRow = $startRow + $tradeRows.Count - 1

# Copy the date/time number format from G4 down the new G column cells
# so the new dates render the same way as the existing ones (style index 1).
$ws.Range("G4").Copy()
$ws.Range("G$startRow`:G$endRow").PasteSpecial(-4122)  # xlPasteFormats

$r = $startRow
foreach ($row in $tradeRows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $r++
}
